# Applies the commit's change: the records in rows 13 and 14 of the
# active sheet are swapped (all "observation" fields - id, sort order,
# red-list category, taxon id, names, author, and coordinates), and the
# now-redundant Starttid/Sluttid ("00:00") cells in columns Z and AB are
# cleared for both rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 13
$row2 = 14

# Columns whose values get swapped between the two rows verbatim.
$swapCols = @("A", "B", "D", "E", "F", "G", "H")

foreach ($col in $swapCols) {
    $cell1 = $ws.Range($col + $row1)
    $cell2 = $ws.Range($col + $row2)
    $v1 = $cell1.Value()
    $v2 = $cell2.Value()
    $cell1.Value = $v2
    $cell2.Value = $v1
}

# Coordinate columns (Ost/Nord) are swapped too, and rounded to whole
# meters in the process.
$coordCols = @("Q", "R")

foreach ($col in $coordCols) {
    $cell1 = $ws.Range($col + $row1)
    $cell2 = $ws.Range($col + $row2)
    $v1 = $cell1.Value()
    $v2 = $cell2.Value()
    $cell1.Value = [Math]::Round([double]$v2)
    $cell2.Value = [Math]::Round([double]$v1)
}

# Starttid / Sluttid columns are no longer populated for either row.
$clearCols = @("Z", "AB")

foreach ($col in $clearCols) {
    $ws.Range($col + $row1).ClearContents()
    $ws.Range($col + $row2).ClearContents()
}
